$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers for ship_A_explore / ship_B_explore
$ws.Range("H1").Value = "ship_A_explore"
$ws.Range("I1").Value = "ship_B_explore"

# Fill explore values of 2 for existing games (rows 2-14, games 1-13)
$ws.Range("H2:I14").Value = 2

# Append new row 15 for game 14
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 3

# Update the active selection to match the final state (I16 selected, matching post-entry cursor)
$ws.Range("I16").Select()
